$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: refresh workflow instance id, UAT gateway id and event date/time
$ws.Range("D2").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.5fc7bb1bb1^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721654277749"
$ws.Range("F2").Value = "22-07-2024:15:17:59"

# Row 3: fill in the previously empty RDA/Regione cell and refresh the IDs/date
$ws.Range("B3").Value = "REGIONE_LAZIO"
$ws.Range("D3").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.308a4c1cc2^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"
$ws.Range("E3").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721654269139"
$ws.Range("F3").Value = "22-07-2024:15:17:51"
